# "finilize sound effect of star and cloud"
#
# The "Music" sheet's row 3 held a placeholder sound-effect entry
# (音效 / www.soundrangers.com / 1.95) that gets finalized: the price/url
# on row 3 are cleared, and row 4 (previously just the placeholder "ding")
# is filled in with the real "sound_star" effect info. The "Music" sheet
# also becomes the active tab/selection instead of "Pistures".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Pistures")
$ws2 = $wb.Worksheets.Item("Music")

# Row 3 on Music keeps its label (音效) but loses the old placeholder
# url/price that used to sit next to it.
$ws2.Range("B3:C3").ClearContents()

# Row 4 becomes the finalized "sound_star" sound effect row.
$ws2.Range("A4").Value = "sound_star"
$ws2.Range("B4").Value = "http://www.2gei.com/sound/class/piano/"
$ws2.Range("C4").Value = "钢琴音符85个_mp3/65!"

# "Pistures" is no longer the selected tab; its selection moves to A4.
$ws1.Range("A4").Select()

# "Music" becomes the active sheet/tab, selection resting on B3.
$ws2.Activate()
$ws2.Range("B3").Select()
